$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the renamed test-case group (SCD0291 -> SCD0018)
$ws.Name = "SCD0018"

# Update TC_ID (column B) for the 4 data rows: "DGS-306" -> "SCD0018-014"
$ws.Range("B2").Value = "SCD0018-014"
$ws.Range("B3").Value = "SCD0018-014"
$ws.Range("B4").Value = "SCD0018-014"
$ws.Range("B5").Value = "SCD0018-014"

# Column B needs to widen a bit to fit the longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.67

# Move the active selection to B6 (scrolls sheet back so column A/B are visible,
# matching the saved view - no more topLeftCell="G1")
$ws.Range("B6").Select()
